$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8; this shifts the old rows 8-13 down to 9-14
# and (like Excel normally does) copies row 7's formatting into the new row 8.
$ws.Rows("8:8").Insert()

# The assignments below are ordered to reproduce the shared-string table order
# of the target file (new strings are appended in first-use order).

# Row 11 (was row 9 before insert): T2, NPN Transistor (BJT)
$ws.Range("A11").Value = "T2"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = " NPN Transistor (BJT)"

# Row 13 (was row 11): R2, R3, R4, R5, R6 / 5 / 200 ohm Resistor
$ws.Range("A13").Value = "R2, R3, R4, R5, R6"
$ws.Range("B13").Value = 5
$ws.Range("C13").Value = "200 Ω Resistor"

# Row 14 (was row 12): U4 / Optocoupler
$ws.Range("A14").Value = "U4"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = " Optocoupler"

# Row 15 (was row 13): S1 / Pushbutton
$ws.Range("A15").Value = "S1"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = " Pushbutton"

# Row 8 (new row): D2 / LASER MODULE
$ws.Range("A8").Value = "D2"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "LASER MODULE"

# Row 9 (was row 8): D2, D4 / 2 / Red LED -> D2 / 1 / Red LED
$ws.Range("A9").Value = "D2"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "Red LED"

# Row 6: component text "Diode" -> "1N4004"
$ws.Range("C6").Value = "1N4004"

# Row 8 link column: laser-module product link, styled like the other hyperlink cells
$ws.Range("D8").Value = "https://www.amazon.com/Youliang-KY-008-Transmitter-Arduino-Raspberry/dp/B07ST98B7S/ref=sr_1_9?keywords=arduino+laser&qid=1675077760&sr=8-9"
$ws.Range("D8").Style = "Hyperlink"

# Row 14 link column: optocoupler product link (plain text, default style)
$ws.Range("D14").Value = "https://www.amazon.com/4N35-Optocouplers-Phototransistor-30V-IC/dp/B0073BOU46"

# Row 15 link column: pushbutton product link (plain text, default style)
$ws.Range("D15").Value = "https://rubikstech.co.ke/index.php?route=product/product&product_id=283"

# Row 10 (was row 9): R1 / 1 kohm Resistor -> BAT1 / 9V Battery + battery product link
$ws.Range("A10").Value = "BAT1"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = " 9V Battery"
$ws.Range("D10").Value = "https://www.amazon.com/Duracell-Coppertop-Alkaline-Batteries-Count/dp/B000K2NW08?th=1"

# Row 11 link column: re-use the existing Digikey link text, styled as a hyperlink
$ws.Range("D11").Value = "https://www.digikey.com/en/products/detail/rochester-electronics-llc/2N3053/12094898"
$ws.Range("D11").Style = "Hyperlink"

# Row 12 (was row 10): R2 / 10 kohm Resistor -> R1 / 1 kohm Resistor (re-uses existing strings)
$ws.Range("A12").Value = "R1"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "1 kΩ Resistor"

# The new plain-text link cells should not inherit a hyperlink look from neighbouring rows
$ws.Range("D10").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D15").ClearFormats()

# Match the saved selection in the target file
$ws.Range("C16").Select()
